# Update Pr(>F) column (F) values per re-run/rounding of PERMANOVA tables.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0.5655
    3  = 0.0088
    4  = 0.5682
    7  = 0.153
    8  = 0.09
    9  = 0.9838
    12 = 0.698
    13 = 0.7278
    14 = 0.5897
    17 = 0.0462
    18 = 0.8251
    19 = 0.3598
    22 = 0.4092
    23 = 0.3778
    24 = 0.1797
    27 = 0.545
    28 = 0.01
    32 = 0.2017
    33 = 0.6376
    34 = 0.9877
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
